$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.378228666666666
$ws.Range("H2").Value = 10.134686
$ws.Range("I2").Value = 0.1977735183221229
$ws.Range("J2").Value = 0.1977735183221229
$ws.Range("M2").Value = 0.01046566666666667
$ws.Range("N2").Value = 0.031397
$ws.Range("O2").Value = 0.007547709568116775
$ws.Range("P2").Value = 0.007547709568116775
$ws.Range("Q2").Value = 0.03535541514911111
$ws.Range("R2").Value = 0.318198736342
$ws.Range("S2").Value = 0.001492737076560005
$ws.Range("T2").Value = 0.001492737076560005
$ws.Range("G3").Value = 3.378228666666666
$ws.Range("H3").Value = 10.134686
$ws.Range("I3").Value = 0.1977735183221229
$ws.Range("J3").Value = 0.1977735183221229
$ws.Range("O3").Value = 0.8438809030711777
$ws.Range("P3").Value = 0.8438809030711776
$ws.Range("Q3").Value = 3.952955448964445
$ws.Range("R3").Value = 35.57659904067999
$ws.Range("S3").Value = 0.1668972952452372
$ws.Range("T3").Value = 0.1668972952452372
$ws.Range("G4").Value = 3.378228666666666
$ws.Range("H4").Value = 10.134686
$ws.Range("I4").Value = 0.1977735183221229
$ws.Range("J4").Value = 0.1977735183221229
$ws.Range("M4").Value = 0.2060093333333333
$ws.Range("N4").Value = 0.618028
$ws.Range("O4").Value = 0.1485713873607056
$ws.Range("P4").Value = 0.1485713873607056
$ws.Range("Q4").Value = 0.6959466354675555
$ws.Range("R4").Value = 6.263519719207999
$ws.Range("S4").Value = 0.02938348600032573
$ws.Range("T4").Value = 0.02938348600032573
$ws.Range("I5").Value = 0.6780480282745078
$ws.Range("J5").Value = 0.6780480282745078
$ws.Range("M5").Value = 0.01046566666666667
$ws.Range("N5").Value = 0.031397
$ws.Range("O5").Value = 0.007547709568116775
$ws.Range("P5").Value = 0.007547709568116775
$ws.Range("Q5").Value = 0.1212127373475556
$ws.Range("R5").Value = 1.090914636128
$ws.Range("S5").Value = 0.005117709590650216
$ws.Range("T5").Value = 0.005117709590650216
$ws.Range("I6").Value = 0.6780480282745078
$ws.Range("J6").Value = 0.6780480282745078
$ws.Range("O6").Value = 0.8438809030711777
$ws.Range("P6").Value = 0.8438809030711776
$ws.Range("S6").Value = 0.5721917824259231
$ws.Range("T6").Value = 0.5721917824259231
$ws.Range("I7").Value = 0.6780480282745078
$ws.Range("J7").Value = 0.6780480282745078
$ws.Range("M7").Value = 0.2060093333333333
$ws.Range("N7").Value = 0.618028
$ws.Range("O7").Value = 0.1485713873607056
$ws.Range("P7").Value = 0.1485713873607056
$ws.Range("Q7").Value = 2.385988012785778
$ws.Range("R7").Value = 21.473892115072
$ws.Range("S7").Value = 0.1007385362579346
$ws.Range("T7").Value = 0.1007385362579346
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.121129333333334
$ws.Range("H8").Value = 6.363388
$ws.Range("I8").Value = 0.1241784534033691
$ws.Range("J8").Value = 0.1241784534033691
$ws.Range("M8").Value = 0.01046566666666667
$ws.Range("N8").Value = 0.031397
$ws.Range("O8").Value = 0.007547709568116775
$ws.Range("P8").Value = 0.007547709568116775
$ws.Range("Q8").Value = 0.02219903255955556
$ws.Range("R8").Value = 0.199791293036
$ws.Range("S8").Value = 0.0009372629009065523
$ws.Range("T8").Value = 0.0009372629009065522
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.121129333333334
$ws.Range("H9").Value = 6.363388
$ws.Range("I9").Value = 0.1241784534033691
$ws.Range("J9").Value = 0.1241784534033691
$ws.Range("O9").Value = 0.8438809030711777
$ws.Range("P9").Value = 0.8438809030711776
$ws.Range("Q9").Value = 2.481989996382223
$ws.Range("R9").Value = 22.33790996744
$ws.Range("S9").Value = 0.1047918254000173
$ws.Range("T9").Value = 0.1047918254000173
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.121129333333334
$ws.Range("H10").Value = 6.363388
$ws.Range("I10").Value = 0.1241784534033691
$ws.Range("J10").Value = 0.1241784534033691
$ws.Range("M10").Value = 0.2060093333333333
$ws.Range("N10").Value = 0.618028
$ws.Range("O10").Value = 0.1485713873607056
$ws.Range("P10").Value = 0.1485713873607056
$ws.Range("Q10").Value = 0.4369724398737779
$ws.Range("R10").Value = 3.932751958864
$ws.Range("S10").Value = 0.01844936510244529
$ws.Range("T10").Value = 0.01844936510244528
